$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: "Cloud 291" (intranet cloud) ---
# Resize/reposition the cloud and tweak its text (two leading spaces before "intranet")
$cloud = $s.Shapes.Item(1)
$cloud.Top = 2917818 / 12700
$cloud.Height = 543310 / 12700
$cloud.TextFrame.TextRange.Text = "  intranet"

# --- Shape 5: "Rectangle 42" ("Core and Middle Tier Modules" -> "Core Modules") ---
$coreModules = $s.Shapes.Item(5)
$trCore = $coreModules.TextFrame.TextRange
# Remove "and Middle Tier " (keeps "Core " and "Modules" as separate runs)
$trCore.Characters(6, 16).Text = ""

# --- Shape 6: "Rectangle 10" ("Voice command" / "input listener" -> "Speech command" / "input listener") ---
$speechCmd = $s.Shapes.Item(6)
$trSpeech = $speechCmd.TextFrame.TextRange
# Merge "input " + "listener" into a single run "input listener"
$trSpeech.Characters(15, 14).Text = "input listener"
# Replace "Voice command" with "Speech command"
$trSpeech.Characters(1, 13).Text = "Speech command"

# --- Shape 8: "Rectangle 69" ("Laser/mouse " + "input " + "listener" -> single run) ---
$laserMouse = $s.Shapes.Item(8)
$trLaser = $laserMouse.TextFrame.TextRange
$trLaser.Characters(1, $trLaser.Text.Length).Text = "Laser/mouse input listener"

# --- Shape 9: "Rectangle 70" ("Tracking postures " + "input " + "listener" -> single run) ---
$tracking = $s.Shapes.Item(9)
$trTrack = $tracking.TextFrame.TextRange
$trTrack.Characters(1, $trTrack.Text.Length).Text = "Tracking postures input listener"
